$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27, shifting existing rows 27-48 down to 28-49
$ws.Rows.Item(27).Insert()

$ws.Cells.Item(27, 1).Value = 5
$ws.Cells.Item(27, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(27, 3).Value = "Maule"
$ws.Cells.Item(27, 4).Value = (Get-Date -Year 2021 -Month 10 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(27, 5).Value = 7
$ws.Cells.Item(27, 6).Value = 100112022
$ws.Cells.Item(27, 7).Value = "Arveja Verde"
$ws.Cells.Item(27, 8).Value = "Sin especificar"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 300
$ws.Cells.Item(27, 11).Value = 17000
$ws.Cells.Item(27, 12).Value = 17000
$ws.Cells.Item(27, 13).Value = 17000
$ws.Cells.Item(27, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(27, 15).Value = "Región del Maule"
$ws.Cells.Item(27, 16).Value = 680
$ws.Cells.Item(27, 17).Value = 25
$ws.Cells.Item(27, 18).Value = "Hortaliza"
